# This script corrects the ordering of a handful of match rows that were
# recorded out of sequence in the "Germany Regionalliga North" results
# table. For each affected pair/trio of rows, the full match record
# (id, div, date, home/away teams, scores, odds, etc. - columns B:AD) is
# swapped/rotated between the rows, while the running index in column A
# is left untouched so the sheet stays sorted by row number.
#
#   Rows   8 <-> 9      (match 6920336 <-> 6920337)
#   Rows  59 <-> 60      (match 6920385 <-> 6920384)
#   Rows 121 -> 122 -> 123 -> 121  (3-way rotation)
#   Rows 215 <-> 216     (match 6919345 <-> 6920542)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 <- original row data from r9
$ws.Cells.Item(8, 2).Value = 6920337
$ws.Cells.Item(8, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(8, 4).Value = [datetime]"2023-07-30 10:00:00"
$ws.Cells.Item(8, 5).Value = "SC SpelleVenhaus"
$ws.Cells.Item(8, 6).Value = "Eimsbutteler TV"
$ws.Cells.Item(8, 7).Value = 3
$ws.Cells.Item(8, 8).Value = 2
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = "H"
$ws.Cells.Item(8, 12).Value = 1.909
$ws.Cells.Item(8, 13).Value = 4
$ws.Cells.Item(8, 14).Value = 3
$ws.Cells.Item(8, 15).Value = 1.7
$ws.Cells.Item(8, 16).Value = 4.333
$ws.Cells.Item(8, 17).Value = 3.5
$ws.Cells.Item(8, 18).Value = -0.75
$ws.Cells.Item(8, 19).Value = 1.95
$ws.Cells.Item(8, 20).Value = 1.9
$ws.Cells.Item(8, 21).Value = 3.25
$ws.Cells.Item(8, 22).Value = 2
$ws.Cells.Item(8, 23).Value = 1.85
$ws.Cells.Item(8, 24).Value = 0.7
$ws.Cells.Item(8, 25).Value = -1
$ws.Cells.Item(8, 26).Value = -1
$ws.Cells.Item(8, 27).Value = 0.475
$ws.Cells.Item(8, 28).Value = -0.5
$ws.Cells.Item(8, 29).Value = 1
$ws.Cells.Item(8, 30).Value = -1

# Row 9 <- original row data from r8
$ws.Cells.Item(9, 2).Value = 6920336
$ws.Cells.Item(9, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(9, 4).Value = [datetime]"2023-07-30 10:00:00"
$ws.Cells.Item(9, 5).Value = "Bremer SV"
$ws.Cells.Item(9, 6).Value = "FC Kilia Kiel"
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 2
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 2
$ws.Cells.Item(9, 11).Value = "D"
$ws.Cells.Item(9, 12).Value = 2.5
$ws.Cells.Item(9, 13).Value = 3.75
$ws.Cells.Item(9, 14).Value = 2.25
$ws.Cells.Item(9, 15).Value = 2.15
$ws.Cells.Item(9, 16).Value = 3.75
$ws.Cells.Item(9, 17).Value = 2.7
$ws.Cells.Item(9, 18).Value = -0.25
$ws.Cells.Item(9, 19).Value = 1.95
$ws.Cells.Item(9, 20).Value = 1.9
$ws.Cells.Item(9, 21).Value = 3
$ws.Cells.Item(9, 22).Value = 1.875
$ws.Cells.Item(9, 23).Value = 1.975
$ws.Cells.Item(9, 24).Value = -1
$ws.Cells.Item(9, 25).Value = 2.75
$ws.Cells.Item(9, 26).Value = -1
$ws.Cells.Item(9, 27).Value = -0.5
$ws.Cells.Item(9, 28).Value = 0.45
$ws.Cells.Item(9, 29).Value = 0.875
$ws.Cells.Item(9, 30).Value = -1

# Row 59 <- original row data from r60
$ws.Cells.Item(59, 2).Value = 6920384
$ws.Cells.Item(59, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(59, 4).Value = [datetime]"2023-09-10 10:00:00"
$ws.Cells.Item(59, 5).Value = "Eintracht Norderstedt"
$ws.Cells.Item(59, 6).Value = "VfB Oldenburg"
$ws.Cells.Item(59, 7).Value = 5
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 11).Value = "H"
$ws.Cells.Item(59, 12).Value = 2.9
$ws.Cells.Item(59, 13).Value = 3.75
$ws.Cells.Item(59, 14).Value = 2
$ws.Cells.Item(59, 15).Value = 3.5
$ws.Cells.Item(59, 16).Value = 4
$ws.Cells.Item(59, 17).Value = 1.75
$ws.Cells.Item(59, 18).Value = 0.75
$ws.Cells.Item(59, 19).Value = 1.8
$ws.Cells.Item(59, 20).Value = 2
$ws.Cells.Item(59, 21).Value = 3.25
$ws.Cells.Item(59, 22).Value = 1.925
$ws.Cells.Item(59, 23).Value = 1.875
$ws.Cells.Item(59, 24).Value = 2.5
$ws.Cells.Item(59, 25).Value = -1
$ws.Cells.Item(59, 26).Value = -1
$ws.Cells.Item(59, 27).Value = 0.8
$ws.Cells.Item(59, 28).Value = -1
$ws.Cells.Item(59, 29).Value = 0.925
$ws.Cells.Item(59, 30).Value = -1

# Row 60 <- original row data from r59
$ws.Cells.Item(60, 2).Value = 6920385
$ws.Cells.Item(60, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(60, 4).Value = [datetime]"2023-09-10 10:00:00"
$ws.Cells.Item(60, 5).Value = "SC SpelleVenhaus"
$ws.Cells.Item(60, 6).Value = "Phonix Lubeck"
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 4
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 2
$ws.Cells.Item(60, 11).Value = "A"
$ws.Cells.Item(60, 12).Value = 9
$ws.Cells.Item(60, 13).Value = 5.5
$ws.Cells.Item(60, 14).Value = 1.222
$ws.Cells.Item(60, 15).Value = 5
$ws.Cells.Item(60, 16).Value = 4
$ws.Cells.Item(60, 17).Value = 1.45
$ws.Cells.Item(60, 18).Value = 1
$ws.Cells.Item(60, 19).Value = 2.05
$ws.Cells.Item(60, 20).Value = 1.8
$ws.Cells.Item(60, 21).Value = 3.25
$ws.Cells.Item(60, 22).Value = 2
$ws.Cells.Item(60, 23).Value = 1.85
$ws.Cells.Item(60, 24).Value = -1
$ws.Cells.Item(60, 25).Value = -1
$ws.Cells.Item(60, 26).Value = 0.45
$ws.Cells.Item(60, 27).Value = -1
$ws.Cells.Item(60, 28).Value = 0.8
$ws.Cells.Item(60, 29).Value = 1
$ws.Cells.Item(60, 30).Value = -1

# Row 121 <- original row data from r123
$ws.Cells.Item(121, 2).Value = 6920445
$ws.Cells.Item(121, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(121, 4).Value = [datetime]"2023-11-04 10:00:00"
$ws.Cells.Item(121, 5).Value = "VfB Oldenburg"
$ws.Cells.Item(121, 6).Value = "Eimsbutteler TV"
$ws.Cells.Item(121, 7).Value = 4
$ws.Cells.Item(121, 8).Value = 1
$ws.Cells.Item(121, 9).Value = 2
$ws.Cells.Item(121, 10).Value = 1
$ws.Cells.Item(121, 11).Value = "H"
$ws.Cells.Item(121, 12).Value = 1.363
$ws.Cells.Item(121, 13).Value = 4.75
$ws.Cells.Item(121, 14).Value = 6
$ws.Cells.Item(121, 15).Value = 1.2
$ws.Cells.Item(121, 16).Value = 5.75
$ws.Cells.Item(121, 17).Value = 9
$ws.Cells.Item(121, 18).Value = -2
$ws.Cells.Item(121, 19).Value = 1.95
$ws.Cells.Item(121, 20).Value = 1.9
$ws.Cells.Item(121, 21).Value = 3.75
$ws.Cells.Item(121, 22).Value = 1.975
$ws.Cells.Item(121, 23).Value = 1.875
$ws.Cells.Item(121, 24).Value = 0.2
$ws.Cells.Item(121, 25).Value = -1
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(121, 27).Value = 0.95
$ws.Cells.Item(121, 28).Value = -1
$ws.Cells.Item(121, 29).Value = 0.9750000000000001
$ws.Cells.Item(121, 30).Value = -1

# Row 122 <- original row data from r121
$ws.Cells.Item(122, 2).Value = 6920446
$ws.Cells.Item(122, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(122, 4).Value = [datetime]"2023-11-04 10:00:00"
$ws.Cells.Item(122, 5).Value = "Bremer SV"
$ws.Cells.Item(122, 6).Value = "St Pauli II"
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 1
$ws.Cells.Item(122, 9).Value = 1
$ws.Cells.Item(122, 10).Value = 1
$ws.Cells.Item(122, 11).Value = "D"
$ws.Cells.Item(122, 12).Value = 3
$ws.Cells.Item(122, 13).Value = 4
$ws.Cells.Item(122, 14).Value = 1.909
$ws.Cells.Item(122, 15).Value = 3.1
$ws.Cells.Item(122, 16).Value = 3.75
$ws.Cells.Item(122, 17).Value = 1.95
$ws.Cells.Item(122, 18).Value = 0.5
$ws.Cells.Item(122, 19).Value = 1.875
$ws.Cells.Item(122, 20).Value = 1.975
$ws.Cells.Item(122, 21).Value = 3.25
$ws.Cells.Item(122, 22).Value = 1.95
$ws.Cells.Item(122, 23).Value = 1.9
$ws.Cells.Item(122, 24).Value = -1
$ws.Cells.Item(122, 25).Value = 2.75
$ws.Cells.Item(122, 26).Value = -1
$ws.Cells.Item(122, 27).Value = 0.875
$ws.Cells.Item(122, 28).Value = -1
$ws.Cells.Item(122, 29).Value = -1
$ws.Cells.Item(122, 30).Value = 0.8999999999999999

# Row 123 <- original row data from r122
$ws.Cells.Item(123, 2).Value = 6920447
$ws.Cells.Item(123, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(123, 4).Value = [datetime]"2023-11-04 10:00:00"
$ws.Cells.Item(123, 5).Value = "Holstein Kiel II"
$ws.Cells.Item(123, 6).Value = "TuS BlauWeiss Lohne"
$ws.Cells.Item(123, 7).Value = 1
$ws.Cells.Item(123, 8).Value = 1
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = "D"
$ws.Cells.Item(123, 12).Value = 1.5
$ws.Cells.Item(123, 13).Value = 4.2
$ws.Cells.Item(123, 14).Value = 4.75
$ws.Cells.Item(123, 15).Value = 1.6
$ws.Cells.Item(123, 16).Value = 4
$ws.Cells.Item(123, 17).Value = 4.2
$ws.Cells.Item(123, 18).Value = -1
$ws.Cells.Item(123, 19).Value = 1.975
$ws.Cells.Item(123, 20).Value = 1.825
$ws.Cells.Item(123, 21).Value = 3.5
$ws.Cells.Item(123, 22).Value = 1.925
$ws.Cells.Item(123, 23).Value = 1.875
$ws.Cells.Item(123, 24).Value = -1
$ws.Cells.Item(123, 25).Value = 3
$ws.Cells.Item(123, 26).Value = -1
$ws.Cells.Item(123, 27).Value = -1
$ws.Cells.Item(123, 28).Value = 0.825
$ws.Cells.Item(123, 29).Value = -1
$ws.Cells.Item(123, 30).Value = 0.875

# Row 215 <- original row data from r216
$ws.Cells.Item(215, 2).Value = 6920542
$ws.Cells.Item(215, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(215, 4).Value = [datetime]"2024-03-24 10:00:00"
$ws.Cells.Item(215, 5).Value = "Holstein Kiel II"
$ws.Cells.Item(215, 6).Value = "VfB Oldenburg"
$ws.Cells.Item(215, 7).Value = 1
$ws.Cells.Item(215, 8).Value = 3
$ws.Cells.Item(215, 9).Value = 1
$ws.Cells.Item(215, 10).Value = 1
$ws.Cells.Item(215, 11).Value = "A"
$ws.Cells.Item(215, 12).Value = 2
$ws.Cells.Item(215, 13).Value = 3.75
$ws.Cells.Item(215, 14).Value = 2.875
$ws.Cells.Item(215, 15).Value = 2.5
$ws.Cells.Item(215, 16).Value = 3.6
$ws.Cells.Item(215, 17).Value = 2.3
$ws.Cells.Item(215, 18).Value = 0
$ws.Cells.Item(215, 19).Value = 1.975
$ws.Cells.Item(215, 20).Value = 1.875
$ws.Cells.Item(215, 21).Value = 3
$ws.Cells.Item(215, 22).Value = 1.875
$ws.Cells.Item(215, 23).Value = 1.975
$ws.Cells.Item(215, 24).Value = -1
$ws.Cells.Item(215, 25).Value = -1
$ws.Cells.Item(215, 26).Value = 1.3
$ws.Cells.Item(215, 27).Value = -1
$ws.Cells.Item(215, 28).Value = 0.875
$ws.Cells.Item(215, 29).Value = 0.875
$ws.Cells.Item(215, 30).Value = -1

# Row 216 <- original row data from r215
$ws.Cells.Item(216, 2).Value = 6919345
$ws.Cells.Item(216, 3).Value = "Germany Regionalliga North"
$ws.Cells.Item(216, 4).Value = [datetime]"2024-03-24 10:00:00"
$ws.Cells.Item(216, 5).Value = "FC Teutonia 05"
$ws.Cells.Item(216, 6).Value = "SSV Jeddeloh"
$ws.Cells.Item(216, 7).Value = 1
$ws.Cells.Item(216, 8).Value = 1
$ws.Cells.Item(216, 9).Value = 1
$ws.Cells.Item(216, 10).Value = 1
$ws.Cells.Item(216, 11).Value = "D"
$ws.Cells.Item(216, 12).Value = 1.6
$ws.Cells.Item(216, 13).Value = 3.75
$ws.Cells.Item(216, 14).Value = 4.5
$ws.Cells.Item(216, 15).Value = 1.7
$ws.Cells.Item(216, 16).Value = 3.5
$ws.Cells.Item(216, 17).Value = 4.2
$ws.Cells.Item(216, 18).Value = -0.75
$ws.Cells.Item(216, 19).Value = 1.95
$ws.Cells.Item(216, 20).Value = 1.85
$ws.Cells.Item(216, 21).Value = 2.75
$ws.Cells.Item(216, 22).Value = 1.775
$ws.Cells.Item(216, 23).Value = 2.025
$ws.Cells.Item(216, 24).Value = -1
$ws.Cells.Item(216, 25).Value = 2.5
$ws.Cells.Item(216, 26).Value = -1
$ws.Cells.Item(216, 27).Value = -1
$ws.Cells.Item(216, 28).Value = 0.8500000000000001
$ws.Cells.Item(216, 29).Value = -1
$ws.Cells.Item(216, 30).Value = 1.025
